# Regenerate the handoff report: record the latest handoff datetime for
# the "05294cb7-..." source file (row 4) on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-03-07 09:31:20"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-03-07 09:31:34"
